# This workbook's data rows (2-15, row 7 excluded) were shuffled: every
# destination row ends up holding the full contents of a different source
# row. Columns C, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX,
# AY are identical across all data rows, so only A, B, D, E, F, G, H, Q, R
# (plus the single AC "Gamla hack" comment that travels from row 13 to row
# 11) actually need to be written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ArtRow {
    param(
        [int]$Row,
        [double]$Id,
        [double]$Taxonsort,
        [string]$Rodlistad,
        [double]$TaxonId,
        [string]$Artnamn,
        [string]$VetNamn,
        [string]$Auktor,
        [double]$Ost,
        [double]$Nord,
        [string]$Kommentar
    )

    $ws.Cells.Item($Row, 1).Value = $Id          # A - Id
    $ws.Cells.Item($Row, 2).Value = $Taxonsort    # B - Taxonsorteringsordning
    $ws.Cells.Item($Row, 4).Value = $Rodlistad    # D - Rödlistade
    $ws.Cells.Item($Row, 5).Value = $TaxonId      # E - TaxonId
    $ws.Cells.Item($Row, 6).Value = $Artnamn      # F - Artnamn
    $ws.Cells.Item($Row, 7).Value = $VetNamn      # G - Vetenskapligt namn
    $ws.Cells.Item($Row, 8).Value = $Auktor       # H - Auktor
    $ws.Cells.Item($Row, 17).Value = $Ost         # Q - Ost
    $ws.Cells.Item($Row, 18).Value = $Nord        # R - Nord

    if ($Kommentar) {
        $ws.Cells.Item($Row, 29).Value = $Kommentar   # AC - Publik kommentar
    } else {
        $ws.Cells.Item($Row, 29).Value = ""
    }
}

Set-ArtRow 2  111523657 89686 "NT" 658    "Rosenticka"        "Rhodofomes roseus"            "(Alb. & Schwein.) Kotl. & Pouzar" 497390.1961838813 6754097.842248607 $null
Set-ArtRow 3  111523724 93881 "LC" 2869   "Bollvitmossa"      "Sphagnum wulfianum"           "Girg."                            497291.3182300103 6754089.649475355 $null
Set-ArtRow 4  111523697 77515 "NT" 6425   "Garnlav"           "Alectoria sarmentosa"         "(Ach.) Ach."                      497380.5053056676 6754165.927741241 $null
Set-ArtRow 5  111523730 89405 "NT" 1202   "Ullticka"          "Phellinidium ferrugineofuscum" "(P.Karst.) Fiasson & Niemelä"    497338.5868253836 6754122.194367126 $null
Set-ArtRow 6  111523712 89405 "NT" 1202   "Ullticka"          "Phellinidium ferrugineofuscum" "(P.Karst.) Fiasson & Niemelä"    497301.0581945881 6754088.183226726 $null
Set-ArtRow 8  111523701 89686 "NT" 658    "Rosenticka"        "Rhodofomes roseus"            "(Alb. & Schwein.) Kotl. & Pouzar" 497367.2942720717 6754083.757028132 $null
Set-ArtRow 9  111523741 89686 "NT" 658    "Rosenticka"        "Rhodofomes roseus"            "(Alb. & Schwein.) Kotl. & Pouzar" 497384.3941364431 6754155.713205664 $null
Set-ArtRow 10 111523728 89686 "NT" 658    "Rosenticka"        "Rhodofomes roseus"            "(Alb. & Schwein.) Kotl. & Pouzar" 497338.5868253836 6754122.194367126 $null
Set-ArtRow 11 111523740 56398 "NT" 100109 "Tretåig hackspett" "Picoides tridactylus"         "(Linnaeus, 1758)"                 497366.3615979423 6754139.679549156 "Gamla hack"
Set-ArtRow 12 111523727 89845 "VU" 1209   "Rynkskinn"         "Phlebia centrifuga"           "P.Karst."                         497338.5868253836 6754122.194367126 $null
Set-ArtRow 13 111523695 5113  "LC" 100526 "Bronshjon"         "Callidium coriaceum"          "Paykull, 1800"                    497354.1644349985 6754111.484663551 $null
Set-ArtRow 14 111523656 89845 "VU" 1209   "Rynkskinn"         "Phlebia centrifuga"           "P.Karst."                         497390.1961838813 6754097.842248607 $null
Set-ArtRow 15 111523683 89845 "VU" 1209   "Rynkskinn"         "Phlebia centrifuga"           "P.Karst."                         497391.6869587752 6754138.20205555  $null
